# Appends a duplicated snapshot of each sheet's existing data rows,
# re-using the same shared-string values (no new unique strings are
# introduced). This mirrors the upstream edit where another reading
# was appended to each of the status/neighbors/links/routes/topology
# sheets in this workbook.

$wb = $excel.ActiveWorkbook

# --- status: duplicate row 2 into row 3 ---
$ws = $wb.Worksheets.Item("status")
$ws.Range("A3:B3").Value2 = $ws.Range("A2:B2").Value2

# --- neighbors: duplicate rows 2-3 into rows 4-5 ---
$ws = $wb.Worksheets.Item("neighbors")
$ws.Range("A4:G5").Value2 = $ws.Range("A2:G3").Value2

# --- links: duplicate rows 2-3 into rows 4-5 ---
$ws = $wb.Worksheets.Item("links")
$ws.Range("A4:G5").Value2 = $ws.Range("A2:G3").Value2

# --- routes: duplicate rows 2-4 into rows 5-7 ---
$ws = $wb.Worksheets.Item("routes")
$ws.Range("A5:G7").Value2 = $ws.Range("A2:G4").Value2

# --- topology: duplicate rows 2-7 into rows 8-13 ---
$ws = $wb.Worksheets.Item("topology")
$ws.Range("A8:G13").Value2 = $ws.Range("A2:G7").Value2
